# electromagnetic break design update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (Загон №1 / A2 = 6)
$ws.Range("B2").Value = 1116.312744140625
$ws.Range("C2").Value = 0.8374
$ws.Range("D2").Value = 0.9006000161170959
$ws.Range("F2").Value = 0.4099999964237213

# Row 3 (A3 = 7)
$ws.Range("B3").Value = 1208.818115234375
$ws.Range("C3").Value = 0.9671
$ws.Range("F3").Value = 0.7092000246047974

# Row 4 (A4 = 8)
$ws.Range("B4").Value = 813.9815063476562
$ws.Range("C4").Value = 0.9509
$ws.Range("D4").Value = 0.9804
$ws.Range("F4").Value = 0.6880999803543091

# Row 5 (A5 = 9)
$ws.Range("B5").Value = 747.32958984375
$ws.Range("C5").Value = 0.7883
$ws.Range("D5").Value = 0.7813
$ws.Range("F5").Value = 0.5238000154495239

# Row 6 (A6 = 10)
$ws.Range("B6").Value = 971.8284912109375
$ws.Range("C6").Value = 0.7689
$ws.Range("D6").Value = 0.7695
$ws.Range("F6").Value = 0.5

# Row 7 (A7 = 11)
$ws.Range("B7").Value = 750.70458984375
$ws.Range("C7").Value = 0.756
$ws.Range("D7").Value = 0.7630000114440918
$ws.Range("F7").Value = 0.5792000293731689

# Row 8 (A8 = 12)
$ws.Range("B8").Value = 766.736572265625
$ws.Range("C8").Value = 0.6870000000000001
$ws.Range("D8").Value = 0.6703
$ws.Range("E8").Value = 0.9742000102996826
$ws.Range("F8").Value = 0.5792000293731689

# Row 9 (A9 = 13)
$ws.Range("B9").Value = 6375.7119140625
$ws.Range("C9").Value = 0.8218
$ws.Range("D9").Value = 0.8110000000000001
$ws.Range("F9").Value = 0.4099999964237213
